$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72; this shifts existing rows 72..103 down to 73..104,
# pushing the former row 103 down to row 104 (matching the target diff).
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new weekly price-report record.
# Columns A, B, C, E, F, G, H, I, K, L, N, O, Q, R repeat the same values used by the
# surrounding rows (same market/category/quality/unit/origin template); only the
# date (D), volume (J), weighted average price (M) and $/Kg price (P) are new.
$row = 72
$ws.Cells.Item($row, 1).Value = 6
$ws.Cells.Item($row, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item($row, 3).Value = "Metropolitana"
$ws.Cells.Item($row, 4).Value = 45135
$ws.Cells.Item($row, 5).Value = 13
$ws.Cells.Item($row, 6).Value = 100112035
$ws.Cells.Item($row, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 430
$ws.Cells.Item($row, 11).Value = 16000
$ws.Cells.Item($row, 12).Value = 17000
$ws.Cells.Item($row, 13).Value = 16465
$ws.Cells.Item($row, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item($row, 15).Value = "Provincia de Quillota"
$ws.Cells.Item($row, 16).Value = 1098
$ws.Cells.Item($row, 17).Value = 15
$ws.Cells.Item($row, 18).Value = "Hortaliza"
